$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.53"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.28%"
$ws.Range("E2").NumberFormat = "General"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.22"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.50%"
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.065"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.07%"
$ws.Range("E4").NumberFormat = "General"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07679"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.95%"
$ws.Range("E5").NumberFormat = "General"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.603"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.09%"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "8.88%"
$ws.Range("E7").NumberFormat = "General"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1247"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.81%"
$ws.Range("E8").NumberFormat = "General"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1855"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.19%"
$ws.Range("E9").NumberFormat = "General"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09191"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.17%"
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04163"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.32%"
$ws.Range("E11").NumberFormat = "General"
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.59%"
$ws.Range("E12").NumberFormat = "General"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001281"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.20%"
$ws.Range("E13").NumberFormat = "General"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "CoinExToken"

$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04160"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "5.04%"
$ws.Range("E14").NumberFormat = "General"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "TigerCash"

$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005757"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.96%"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "UpBots"

$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.007430"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1,888.57%"
$ws.Range("E16").NumberFormat = "General"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "LEO"

$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.328"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.36%"
$ws.Range("E17").NumberFormat = "General"
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "GateToken"

$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.413"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.18%"
$ws.Range("E18").NumberFormat = "General"
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "BTSEToken"

$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.356"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.89%"
$ws.Range("E19").NumberFormat = "General"
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = "BitpandaEcosystemToken"

$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3353"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.85%"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").Value = "MCDex"

$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.416"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.17%"
$ws.Range("E21").NumberFormat = "General"
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").Value = "ProBitToken"

$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1399"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.01%"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").Value = "ZBToken"

$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.3198"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "13.83%"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001280"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.64%"
$ws.Range("E24").NumberFormat = "General"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "18.19%"
$ws.Range("E25").NumberFormat = "General"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.50%"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02448"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "5.31%"
$ws.Range("E38").NumberFormat = "General"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05280"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.16%"
$ws.Range("E39").NumberFormat = "General"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005878"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.83%"
$ws.Range("E40").NumberFormat = "General"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007652"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.16%"
$ws.Range("E41").NumberFormat = "General"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1346"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.32%"
$ws.Range("E42").NumberFormat = "General"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007338"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.82%"
$ws.Range("E43").NumberFormat = "General"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007457"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.15%"
$ws.Range("E44").NumberFormat = "General"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3020"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.51%"
$ws.Range("E45").NumberFormat = "General"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006682"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.51%"
$ws.Range("E46").NumberFormat = "General"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.51%"
$ws.Range("E47").NumberFormat = "General"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03994"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-13.78%"
$ws.Range("E48").NumberFormat = "General"
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.37%"
$ws.Range("E49").NumberFormat = "General"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.51%"
$ws.Range("E50").NumberFormat = "General"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.51%"
$ws.Range("E51").NumberFormat = "General"
$ws.Range("E51").Style = "Normal"
